# Update hookworm coverage scenario 2 workbook:
#  - Expand the yearly timeline (2018-2040) on both sheets into a half-year
#    timeline (2018, 2018.5, 2019, 2019.5, ..., 2040).
#  - On "Platform Coverage": from 2026 onward, MDA coverage rows (0.8 / 0.5 / 0.5)
#    are applied continuously every half year instead of every other year.
#  - On "MarketShare": the "Old Product B (SOC)" market share of 1 is extended
#    across the full (now half-yearly) timeline through 2040.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Platform Coverage"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Platform Coverage")

# First column of the year timeline (column H = 8) through the new last
# column (column AZ = 52), one column per half-year step.
$startCol = 8
$endCol = 52

# --- Row 1: half-year header values ------------------------------------
$year = 2018.0
for ($col = $startCol; $col -le $endCol; $col++) {
    $ws1.Cells.Item(1, $col).Value = $year
    $year = $year + 0.5
}

# --- Row 2: 0.6 MDA coverage (ages 5-15), annual cadence 2018-2025 -------
# Columns H, J, L, N, P, R, T, V correspond (in the new half-year grid) to
# 2018, 2019, 2020, 2021, 2022, 2023, 2024, 2025.
$row2cols = @(8, 10, 12, 14, 16, 18, 20, 22)
foreach ($c in $row2cols) {
    $ws1.Cells.Item(2, $c).Value = 0.6
}

# --- Rows 3-5: continuous half-year coverage from 2026 (column X = 24) ---
# through 2040 (column AZ = 52). The old biennial columns (P, R, T, V =
# 16, 18, 20, 22) that used to carry this value are cleared since in the
# new half-year grid they now represent different (earlier) years.
$oldBiennialCols = @(16, 18, 20, 22)
foreach ($c in $oldBiennialCols) {
    $ws1.Cells.Item(3, $c).ClearContents()
    $ws1.Cells.Item(4, $c).ClearContents()
    $ws1.Cells.Item(5, $c).ClearContents()
}

$contStart = 24  # column X -> 2026
$contEnd = 52    # column AZ -> 2040
for ($c = $contStart; $c -le $contEnd; $c++) {
    $ws1.Cells.Item(3, $c).Value = 0.8
    $ws1.Cells.Item(4, $c).Value = 0.5
    $ws1.Cells.Item(5, $c).Value = 0.5
}

# ---------------------------------------------------------------------------
# Sheet 2: "MarketShare"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("MarketShare")

# Year timeline starts at column D (4) through column AV (48).
$startCol2 = 4
$endCol2 = 48

# --- Row 1: half-year header values ------------------------------------
$year2 = 2018.0
for ($col = $startCol2; $col -le $endCol2; $col++) {
    $ws2.Cells.Item(1, $col).Value = $year2
    $year2 = $year2 + 0.5
}

# --- Row 3: "Old Product B (SOC)" market share = 1 across full timeline --
for ($col = $startCol2; $col -le $endCol2; $col++) {
    $ws2.Cells.Item(3, $col).Value = 1
}

# ---------------------------------------------------------------------------
# Restore the final on-screen selection for each sheet (matches the saved
# workbook state) and leave "Platform Coverage" as the active sheet.
# ---------------------------------------------------------------------------
$ws1.Select()
$ws1.Range("R19").Select()
$ws2.Select()
$ws2.Range("AW3").Select()
$ws1.Select()
